$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1, styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2:F11 with the time_taken values, stored as text
$timeTaken = @(
    "2021-10-05 13:41:56.709831",
    "2021-10-05 13:41:56.709844",
    "2021-10-05 13:41:56.709848",
    "2021-10-05 13:41:56.709851",
    "2021-10-05 13:41:56.709854",
    "2021-10-05 13:41:56.709857",
    "2021-10-05 13:41:56.709860",
    "2021-10-05 13:41:56.709864",
    "2021-10-05 13:41:56.709867",
    "2021-10-05 13:41:56.709870"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timeTaken[$i]
}
